$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.153.07"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "1.785.54"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "2.043.59"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").Value = "1.783.41"
$ws.Range("E14").Value = "  +1.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("D16").Value = "34.093.64"

$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.21%  "

$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("E22").Value = "  +3.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.22%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.71%  "

$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0519"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").Value = "1.449.49"
$ws.Range("E35").Value = "  +5.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.656"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.05%  "

$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.922"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("E48").Value = "  -4.81%  "

$ws.Range("D49").Value = "1.945.19"
$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
